$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("F15").Value = "C3225X7R1A226K230AC"
$ws.Range("D15").Value = "'1210'"
$ws.Range("G15").Value = "'tdk"

$ws.Range("F15").Font.Size = 9
$ws.Range("F15").Font.Name = "Arial"
$ws.Range("F15").Font.Color = 4473924
